$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 2144.3333
$ws.Range("J29").Value = 4214.5
$ws.Range("L29").Value = 12643.5
$ws.Range("N29").Value = -13205.5
$ws.Range("H38").Value = 3915.6
$ws.Range("I38").Value = 1651.2
$ws.Range("J38").Value = 6180
$ws.Range("K38").Value = 4953.6
$ws.Range("L38").Value = 18540
$ws.Range("M38").Value = -4581.6
$ws.Range("N38").Value = -19284
$ws.Range("H40").Value = 4898.6787
$ws.Range("I40").Value = 3422.1
$ws.Range("J40").Value = 8590.125
$ws.Range("K40").Value = 3422.1
$ws.Range("L40").Value = 8590.125
$ws.Range("M40").Value = -3247.1
$ws.Range("N40").Value = -8940.125
$ws.Range("H41").Value = 1189.0714
$ws.Range("I41").Value = 731.5
$ws.Range("J41").Value = 1799.1666
$ws.Range("K41").Value = 731.5
$ws.Range("L41").Value = 1799.1666
$ws.Range("M41").Value = -291.5
$ws.Range("N41").Value = -2679.1666
$ws.Range("H58").Value = 2117.4443
$ws.Range("J58").Value = 3000
$ws.Range("L58").Value = 9000
$ws.Range("N58").Value = -9300
$ws.Range("H62").Value = 8101.5713
$ws.Range("I62").Value = 6397
$ws.Range("K62").Value = 6397
$ws.Range("M62").Value = -5773
$ws.Range("H65").Value = 8101.5713
$ws.Range("I65").Value = 6397
$ws.Range("K65").Value = 31985
$ws.Range("M65").Value = -28865
$ws.Range("H86").Value = 4133.8335
$ws.Range("I86").Value = 3160.6
$ws.Range("J86").Value = 9000
$ws.Range("K86").Value = 3160.6
$ws.Range("L86").Value = 9000
$ws.Range("M86").Value = -2037.6
$ws.Range("N86").Value = -11246
$ws.Range("H89").Value = 4133.8335
$ws.Range("I89").Value = 3160.6
$ws.Range("J89").Value = 9000
$ws.Range("K89").Value = 15803
$ws.Range("L89").Value = 45000
$ws.Range("M89").Value = -10187
$ws.Range("N89").Value = -56232
$ws.Range("H101").Value = 168
$ws.Range("I101").Value = 168
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 504
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = 1118
$ws.Range("N101").Value = ""
$ws.Range("H106").Value = 10000
$ws.Range("I106").Value = 10000
$ws.Range("K106").Value = 10000
$ws.Range("M106").Value = -9369
$ws.Range("H115").Value = 500
$ws.Range("I115").Value = 500
$ws.Range("K115").Value = 1500
$ws.Range("M115").Value = 67
$ws.Range("H137").Value = 4017.2354
$ws.Range("I137").Value = 2099.8333
$ws.Range("K137").Value = 6299.499899999999
$ws.Range("M137").Value = -3749.499899999999
$ws.Range("H138").Value = 3387.5625
$ws.Range("I138").Value = 3199.4285
$ws.Range("J138").Value = 3533.889
$ws.Range("K138").Value = 9598.2855
$ws.Range("L138").Value = 10601.667
$ws.Range("M138").Value = -4458.2855
$ws.Range("N138").Value = -20881.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1020.82355
$ws.Range("I2").Value = 918.7857
$ws.Range("K2").Value = 918.7857
$ws.Range("M2").Value = -805.7857
$ws.Range("H45").Value = 2486.7693
$ws.Range("I45").Value = 2104.4443
$ws.Range("K45").Value = 2104.4443
$ws.Range("M45").Value = -1727.4443
$ws.Range("H61").Value = 4585.1177
$ws.Range("I61").Value = 2995.182
$ws.Range("K61").Value = 2995.182
$ws.Range("M61").Value = -2783.182
$ws.Range("H116").Value = 1020.82355
$ws.Range("I116").Value = 918.7857
$ws.Range("K116").Value = 918.7857
$ws.Range("M116").Value = 1375.2143
$ws.Range("H130").Value = 20428
$ws.Range("J130").Value = 20428
$ws.Range("L130").Value = 20428
$ws.Range("N130").Value = -30468
$ws.Range("H132").Value = 5409.579
$ws.Range("I132").Value = 4252.4614
$ws.Range("K132").Value = 12757.3842
$ws.Range("M132").Value = -10227.3842
$ws.Range("H136").Value = 4585.1177
$ws.Range("I136").Value = 2995.182
$ws.Range("K136").Value = 8985.545999999998
$ws.Range("M136").Value = -6435.545999999998
$ws.Range("H139").Value = 23999
$ws.Range("J139").Value = 23999
$ws.Range("L139").Value = 23999
$ws.Range("N139").Value = -34279

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1020.82355
$ws.Range("I3").Value = 918.7857
$ws.Range("K3").Value = 918.7857
$ws.Range("M3").Value = -804.7857
$ws.Range("H22").Value = 1261.5454
$ws.Range("I22").Value = 1378.7
$ws.Range("J22").Value = 90
$ws.Range("K22").Value = 1378.7
$ws.Range("L22").Value = 90
$ws.Range("M22").Value = -1205.7
$ws.Range("N22").Value = -436
$ws.Range("H86").Value = 6458.65
$ws.Range("I86").Value = 4925.091
$ws.Range("J86").Value = 8333
$ws.Range("K86").Value = 4925.091
$ws.Range("L86").Value = 8333
$ws.Range("M86").Value = -3802.091
$ws.Range("N86").Value = -10579
$ws.Range("H89").Value = 6458.65
$ws.Range("I89").Value = 4925.091
$ws.Range("J89").Value = 8333
$ws.Range("K89").Value = 24625.455
$ws.Range("L89").Value = 41665
$ws.Range("M89").Value = -19009.455
$ws.Range("N89").Value = -52897
$ws.Range("H102").Value = 14999.5
$ws.Range("I102").Value = 14999.5
$ws.Range("K102").Value = 14999.5
$ws.Range("M102").Value = -11754.5
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6511.8887
$ws.Range("I31").Value = 1638
$ws.Range("J31").Value = 8564.053
$ws.Range("K31").Value = 1638
$ws.Range("L31").Value = 8564.053
$ws.Range("M31").Value = -1343
$ws.Range("N31").Value = -9154.053
$ws.Range("H34").Value = 6511.8887
$ws.Range("I34").Value = 1638
$ws.Range("J34").Value = 8564.053
$ws.Range("K34").Value = 1638
$ws.Range("L34").Value = 8564.053
$ws.Range("M34").Value = -1436
$ws.Range("N34").Value = -8968.053
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").Value = ""
$ws.Range("H86").Value = 5169.3335
$ws.Range("I86").Value = 2750
$ws.Range("K86").Value = 2750
$ws.Range("M86").Value = -1627
$ws.Range("H89").Value = 5169.3335
$ws.Range("I89").Value = 2750
$ws.Range("K89").Value = 13750
$ws.Range("M89").Value = -8134
$ws.Range("H105").Value = 2231.3572
$ws.Range("I105").Value = 2159.889
$ws.Range("K105").Value = 2159.889
$ws.Range("M105").Value = -412.8890000000001
$ws.Range("H106").Value = 28556.666
$ws.Range("J106").Value = 28556.666
$ws.Range("L106").Value = 28556.666
$ws.Range("N106").Value = -31080.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 815.82355
$ws.Range("I5").Value = 771.8889
$ws.Range("J5").Value = 865.25
$ws.Range("K5").Value = 2315.6667
$ws.Range("L5").Value = 2595.75
$ws.Range("M5").Value = -2203.6667
$ws.Range("N5").Value = -2819.75
$ws.Range("H39").Value = 6344.0586
$ws.Range("J39").Value = 7147.933
$ws.Range("L39").Value = 21443.799
$ws.Range("N39").Value = -22031.799
$ws.Range("H92").Value = 659.6667
$ws.Range("I92").Value = 650.5
$ws.Range("J92").Value = 664.25
$ws.Range("K92").Value = 1951.5
$ws.Range("L92").Value = 1992.75
$ws.Range("M92").Value = -703.5
$ws.Range("N92").Value = -4488.75
$ws.Range("H135").Value = 815.82355
$ws.Range("I135").Value = 771.8889
$ws.Range("J135").Value = 865.25
$ws.Range("K135").Value = 6947.0001
$ws.Range("L135").Value = 7787.25
$ws.Range("M135").Value = -4412.0001
$ws.Range("N135").Value = -12857.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 59745.668
$ws.Range("J140").Value = 69770
$ws.Range("L140").Value = 69770
$ws.Range("N140").Value = -80130

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3000
$ws.Range("J132").Value = 3000
$ws.Range("L132").Value = 9000
$ws.Range("N132").Value = -14060
$ws.Range("H136").Value = 4663.154
$ws.Range("I136").Value = 3311
$ws.Range("K136").Value = 9933
$ws.Range("M136").Value = -7383

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 37898.4
$ws.Range("J103").Value = 37898.4
$ws.Range("L103").Value = 37898.4
$ws.Range("N103").Value = -40242.4
$ws.Range("H113").Value = 605.94116
$ws.Range("I113").Value = 516.9
$ws.Range("K113").Value = 1550.7
$ws.Range("M113").Value = 619.3000000000002
$ws.Range("H122").Value = 5461
$ws.Range("I122").Value = 4444
$ws.Range("K122").Value = 13332
$ws.Range("M122").Value = -10882
$ws.Range("H132").Value = 1935.7
$ws.Range("I132").Value = 1862.5
$ws.Range("K132").Value = 5587.5
$ws.Range("M132").Value = -3057.5
$ws.Range("H133").Value = 100000
$ws.Range("J133").Value = 100000
$ws.Range("L133").Value = 100000
$ws.Range("N133").Value = -110120
$ws.Range("H136").Value = 2887.8
$ws.Range("I136").Value = 2066.389
$ws.Range("K136").Value = 6199.167
$ws.Range("M136").Value = -3649.167
